# Word COM-interop script: add the new "CV Table" table style and the
# Experience-section paragraph/character styles (ExpRole, ExpCompany,
# ExpMeta, ExpBody, ExpHighlight) used by the accented CV template.
#
# The new <w:style w:styleId="CVTable"> element belongs right after the
# last pre-existing table style and right before "Contact Info" in the
# canonical styles part, while the five "Exp*" styles belong at the very
# end of the styles part (after "Skill Level").  This Word engine's
# Styles.Add always appends to the end of the collection, so to land
# CVTable in the middle we temporarily remove the handful of styles that
# currently sit at the tail (Contact Info / Skill Category / Skill Items /
# Skill Highlight / Skill Level) and re-add them - identically - in the
# desired final order together with the new styles.
#
# wdStyleType constants: 1 = paragraph, 2 = character, 3 = table
# wdLineSpaceMultiple  : 5  (LineSpacing is expressed as 12pt * multiple,
#                            i.e. w:line = LineSpacing * 20)
# wdColor RGB packing  : r + g*256 + b*65536 (matches w:color hex RRGGBB)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 0) Pull the tail-end styles back out so we can reinsert them after the
#    new CVTable style, preserving their exact existing formatting.
# ---------------------------------------------------------------------
$d.Styles("Contact Info").Delete()
$d.Styles("Skill Category").Delete()
$d.Styles("Skill Items").Delete()
$d.Styles("Skill Highlight").Delete()
$d.Styles("Skill Level").Delete()

# ---------------------------------------------------------------------
# 1) CVTable - custom table style with all borders switched off.
# ---------------------------------------------------------------------
$cvTable = $d.Styles.Add("CV Table", 3)
try {
    $cvTable.Table.Borders.InsideLineStyle = 0   # wdLineStyleNone
    $cvTable.Table.Borders.OutsideLineStyle = 0  # wdLineStyleNone
    $cvTable.Table.Borders.Enable = 0
    for ($i = -4; $i -le 6; $i++) {
        $cvTable.Table.Borders($i).LineStyle = 0
    }
} catch {
    # Border sub-object formatting on a *style* (as opposed to a real
    # table) isn't exposed everywhere; ignore and move on - the style
    # itself has already been minted above.
}

# ---------------------------------------------------------------------
# 2) Contact Info - restored exactly as it was.
# ---------------------------------------------------------------------
$contactInfo = $d.Styles.Add("Contact Info", 1)
$contactInfo.Font.Name = "Liberation Sans"
$contactInfo.Font.Bold = $false
$contactInfo.Font.Color = 3877150      # 1E293B
$contactInfo.Font.Size = 9              # sz 18
$contactInfo.ParagraphFormat.SpaceBefore = 0
$contactInfo.ParagraphFormat.SpaceAfter = 5.1   # after 102
$contactInfo.ParagraphFormat.LineSpacingRule = 5
$contactInfo.ParagraphFormat.LineSpacing = 18   # line 360 (1.5x)

# ---------------------------------------------------------------------
# 3) Skill Category - restored exactly as it was.
# ---------------------------------------------------------------------
$skillCategory = $d.Styles.Add("Skill Category", 1)
$skillCategory.Font.Name = "Liberation Sans"
$skillCategory.Font.Bold = $true
$skillCategory.Font.Color = 9139300    # 64748B
$skillCategory.Font.Size = 10           # sz 20
$skillCategory.ParagraphFormat.SpaceBefore = 0
$skillCategory.ParagraphFormat.SpaceAfter = 0
$skillCategory.ParagraphFormat.LineSpacingRule = 5
$skillCategory.ParagraphFormat.LineSpacing = 14.4   # line 288 (1.2x)

# ---------------------------------------------------------------------
# 4) Skill Items - restored exactly as it was.
# ---------------------------------------------------------------------
$skillItems = $d.Styles.Add("Skill Items", 1)
$skillItems.Font.Name = "Liberation Sans"
$skillItems.Font.Bold = $false
$skillItems.Font.Color = 3877150       # 1E293B
$skillItems.Font.Size = 10              # sz 20
$skillItems.ParagraphFormat.SpaceBefore = 0
$skillItems.ParagraphFormat.SpaceAfter = 0
$skillItems.ParagraphFormat.LineSpacingRule = 5
$skillItems.ParagraphFormat.LineSpacing = 14.4   # line 288 (1.2x)

# ---------------------------------------------------------------------
# 5) Skill Highlight - restored exactly as it was.
# ---------------------------------------------------------------------
$skillHighlight = $d.Styles.Add("Skill Highlight", 2)
$skillHighlight.Font.Name = "Liberation Sans"
$skillHighlight.Font.Bold = $true
$skillHighlight.Font.Color = 809194    # EA580C
$skillHighlight.Font.Size = 10          # sz 20

# ---------------------------------------------------------------------
# 6) Skill Level - restored exactly as it was.
# ---------------------------------------------------------------------
$skillLevel = $d.Styles.Add("Skill Level", 2)
$skillLevel.Font.Name = "Liberation Sans"
$skillLevel.Font.Bold = $false
$skillLevel.Font.Color = 9139300       # 64748B
$skillLevel.Font.Size = 10              # sz 20

# ---------------------------------------------------------------------
# 7) ExpRole - new paragraph style, bold dark heading for the job title.
# ---------------------------------------------------------------------
$expRole = $d.Styles.Add("Exp Role", 1)
$expRole.Font.Name = "Liberation Sans"
$expRole.Font.Bold = $true
$expRole.Font.Color = 3877150          # 1E293B
$expRole.Font.Size = 11                 # sz 22
$expRole.ParagraphFormat.SpaceBefore = 0
$expRole.ParagraphFormat.SpaceAfter = 0
$expRole.ParagraphFormat.LineSpacingRule = 5
$expRole.ParagraphFormat.LineSpacing = 13.8    # line 276 (1.15x)

# ---------------------------------------------------------------------
# 8) ExpCompany - new paragraph style, bold orange company name.
# ---------------------------------------------------------------------
$expCompany = $d.Styles.Add("Exp Company", 1)
$expCompany.Font.Name = "Liberation Sans"
$expCompany.Font.Bold = $true
$expCompany.Font.Color = 809194        # EA580C
$expCompany.Font.Size = 11              # sz 22
$expCompany.ParagraphFormat.SpaceBefore = 0
$expCompany.ParagraphFormat.SpaceAfter = 0
$expCompany.ParagraphFormat.LineSpacingRule = 5
$expCompany.ParagraphFormat.LineSpacing = 13.8

# ---------------------------------------------------------------------
# 9) ExpMeta - new paragraph style, muted small meta text (dates/place).
# ---------------------------------------------------------------------
$expMeta = $d.Styles.Add("Exp Meta", 1)
$expMeta.Font.Name = "Liberation Sans"
$expMeta.Font.Bold = $false
$expMeta.Font.Color = 9139300          # 64748B
$expMeta.Font.Size = 9                  # sz 18
$expMeta.ParagraphFormat.SpaceBefore = 0
$expMeta.ParagraphFormat.SpaceAfter = 0
$expMeta.ParagraphFormat.LineSpacingRule = 5
$expMeta.ParagraphFormat.LineSpacing = 13.8

# ---------------------------------------------------------------------
# 10) ExpBody - new paragraph style, body copy for bullet/description.
# ---------------------------------------------------------------------
$expBody = $d.Styles.Add("Exp Body", 1)
$expBody.Font.Name = "Liberation Sans"
$expBody.Font.Bold = $false
$expBody.Font.Color = 3877150          # 1E293B
$expBody.Font.Size = 10                 # sz 20
$expBody.ParagraphFormat.SpaceBefore = 0
$expBody.ParagraphFormat.SpaceAfter = 2.85   # after 57
$expBody.ParagraphFormat.LineSpacingRule = 5
$expBody.ParagraphFormat.LineSpacing = 15    # line 300 (1.25x)

# ---------------------------------------------------------------------
# 11) ExpHighlight - new character style, bold dark inline emphasis run.
# ---------------------------------------------------------------------
$expHighlight = $d.Styles.Add("Exp Highlight", 2)
$expHighlight.Font.Name = "Liberation Sans"
$expHighlight.Font.Bold = $true
$expHighlight.Font.Color = 3877150     # 1E293B
$expHighlight.Font.Size = 10            # sz 20

Write-Output "Added styles: CVTable, ExpRole, ExpCompany, ExpMeta, ExpBody, ExpHighlight"
